$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell-level value changes (column C impute/un-impute) above row 26 ---
$ws.Range("C2").Value = 14.9
$ws.Range("C6").ClearContents()
$ws.Range("C12").Value = 12.5
$ws.Range("C14").ClearContents()
$ws.Range("C20").Value = 12.5
$ws.Range("C21").Value = 12.7
$ws.Range("C23").ClearContents()
$ws.Range("C24").ClearContents()

# --- Row removals: "RM 232" row and "SC 92" row were dropped entirely ---
# Row 26 is "RM 232" in the original layout.
$ws.Rows(26).Delete()
# After that deletion, the "SC 92" row (originally row 28) is now row 27.
$ws.Rows(27).Delete()

# --- Remaining rows (now 26-33) need some B/C values adjusted ---
$ws.Range("B26").Value = -20.2      # SC 5
$ws.Range("B27").ClearContents()    # SC 101
$ws.Range("B28").ClearContents()    # SC 105
$ws.Range("B29").Value = -19.5      # SC 119
$ws.Range("B30").Value = -19.7      # SC 120
$ws.Range("B31").ClearContents()    # SC 132
$ws.Range("C31").Value = 15.3       # SC 132
$ws.Range("B32").ClearContents()    # SC 193
$ws.Range("C33").Value = 10.4       # SC 232
